# Lab3 results.xlsx update:
#  - Add a "time(s)" (speedup) column C to both data tables
#  - Re-point the two bar charts at the new column C data
#  - Add Max/Min scale to the "With load balancing" chart's value axis
#  - Add a small "Series / Speedup / SpeedupLB" summary block
#  - Center-align + merge the two title rows across A:C
#  - Auto-size column A, update the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header for the new "time(s)" column, table 1 (rows 2-6)
# ---------------------------------------------------------------------
$ws.Range("C2").Value = "time(s)"
$ws.Range("C3").Value = 1.37002
$ws.Range("C4").Value = 3.5270100000000002
$ws.Range("C5").Value = 5.2553900000000002
$ws.Range("C6").Value = 6.4890999999999996

# ---------------------------------------------------------------------
# 2. Header for the new "time(s)" column, table 2 (rows 9-13)
# ---------------------------------------------------------------------
$ws.Range("C9").Value = "time(s)"
$ws.Range("C10").Value = 5.4816900000000004
$ws.Range("C11").Value = 5.5274999999999999
$ws.Range("C12").Value = 5.5132599999999998
$ws.Range("C13").Value = 5.5763100000000003

# ---------------------------------------------------------------------
# 3. New summary block (rows 19-21)
# ---------------------------------------------------------------------
$ws.Range("A19").Value = "Series"
$ws.Range("B19").Value = 10.191000000000001

$ws.Range("A20").Value = "Speedup"
$ws.Range("B20").Formula = "=B19/C6"

$ws.Range("A21").Value = "SpeedupLB"
$ws.Range("B21").Formula = "=B19/C13"

# ---------------------------------------------------------------------
# 4. Re-point the two charts at the new "time(s)" column
# ---------------------------------------------------------------------
$chart1 = $ws.ChartObjects().Item(1).Chart
$chart1.SeriesCollection(1).Values = "Sheet1!`$C`$3:`$C`$6"

$chart2 = $ws.ChartObjects().Item(2).Chart
$chart2.SeriesCollection(1).Values = "Sheet1!`$C`$10:`$C`$13"

# Give the "With load balancing" chart an explicit value-axis scale
$valAx2 = $chart2.Axes(2)
$valAx2.MinimumScale = 0
$valAx2.MaximumScale = 7

# ---------------------------------------------------------------------
# 5. Title rows: merge A:C and center them
# ---------------------------------------------------------------------
$ws.Range("A1:C1").HorizontalAlignment = -4108
[void]$ws.Range("A1:C1").Merge()

$ws.Range("A8:C8").HorizontalAlignment = -4108
[void]$ws.Range("A8:C8").Merge()

# ---------------------------------------------------------------------
# 6. Column A width + selection
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 10.85546875
[void]$ws.Range("F27").Select()
